# Generate Report for Handback
# Updates the handback-status report with newly generated xliff handoff/handback
# timestamps for the second file (9f945c4e-...) in each sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-16 00:44:10"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-16 00:44:02"
$wsZhCn.Range("K3").Value = "2016-08-16 00:44:27"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-16 00:44:10"
$wsDeDe.Range("K3").Value = "2016-08-16 00:44:34"
